$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's gold-price row (row 71), mirroring the format of the
# preceding row (row 70: date in column A, wrapped price text in column B).
$ws.Range("A70:B70").Copy()
$ws.Range("A71:B71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A71").Value = "26-11-2025"
$ws.Range("B71").Value = "The price of gold in India today is ₹12,791 per gram for 24 karat gold, ₹11,725 per gram for 22 karat gold and ₹9,593 per gram for 18 karat gold (also called 999 gold)."
